$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.706.09'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.05%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.685.96'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +2.53%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '220.39'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.17%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.521'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.00'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '30.57'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +3.85%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.265'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +1.32%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +1.71%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.931.48'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +2.77%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '10.55'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  +11.16%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.623'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +8.54%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '1.686.50'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +2.50%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.98'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.48%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '30.737.70'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '66.43'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +2.24%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '245.76'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.93%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0714'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.38%  '
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.11%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.21'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +1.61%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.27'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +1.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.15'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -0.06%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.28'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.27%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.85'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.46%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.69'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +0.05%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0499'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +0.81%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.41%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.48'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +1.46%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.510.21'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +4.69%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.28'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +1.73%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.73'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +3.53%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -1.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '83.49'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +6.11%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +3.91%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.73'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -4.87%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +4.36%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.40%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.846'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.99'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.63%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +0.52%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.29%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.00'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '51.12'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -8.08%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.824.24'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +2.07%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '5.42'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +0.52%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '94.37'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +4.36%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0₆0113'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -0.81%  '
